$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab (Book1 -> SAT_Test_1)
$ws.Name = "SAT_Test_1"

# Insert a new column before column A; this shifts the existing
# Question/Answer columns (old A:E) to B:F, carrying their values,
# shared-string refs, styles (s="1") and the custom column width along.
$ws.Range("A1").EntireColumn.Insert()

# Populate the new column A with the question "Type"
$ws.Range("A1").Value = "Type"
$ws.Range("A2").Value = "Math"
$ws.Range("A3").Value = "Math"
$ws.Range("A4").Value = "Math"
$ws.Range("A5").Value = "Math"
$ws.Range("A6").Value = "Math"

# New row 6: an image-based question
$ws.Range("B6").Value = "Which of the following graphs best shows a strong negative association between d and t?"

# New columns G (QuestionImage) and H (Answer Image)
$ws.Range("G1").Value = "QuestionImage"
$ws.Range("H1").Value = "Answer Image"

$ws.Range("G2").Value = "F"
$ws.Range("H2").Value = "F"
$ws.Range("G3").Value = "F"
$ws.Range("H3").Value = "F"
$ws.Range("G4").Value = "F"
$ws.Range("H4").Value = "F"
$ws.Range("G5").Value = "F"
$ws.Range("H5").Value = "F"
$ws.Range("G6").Value = "F"
$ws.Range("H6").Value = "T"

# Apply the same formatting (wrap text / top alignment, text number format)
# used by the existing question/answer cells to every newly written cell,
# by copying format from an already-styled cell (B1, which carries s="1").
# Target only the cells that actually hold content in the new layout so we
# don't leave stray styled-but-empty cells behind (row 6 only has B/G/H).
$ws.Range("B1").Copy()
$ws.Range("G1:H5").PasteSpecial(-4122)
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("G6:H6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row heights: header row grows to fit the wrapped "Type"/"Question" text,
# and the new image-question row needs room too.
$ws.Rows.Item(1).RowHeight = 30
$ws.Rows.Item(6).RowHeight = 45

# Selection state
$ws.Range("J4").Select()
